# Fill every blank cell in the "Kotoba/bab2" vocabulary sheet with "-".
# The original sheet had a bunch of empty-but-styled cells (no <v>, no t="s")
# scattered across columns B:E, rows 2:50. The edit fills each one with the
# literal text "-", including one brand new cell (B50) that did not exist
# in the sheet before.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$blankCells = @(
    "C2", "C3", "C4", "C5", "C6", "C7", "E8", "E9", "E10", "E11",
    "C12", "E12", "E13", "E14", "C15", "E15", "C16", "E16", "E17", "C18",
    "E18", "C19", "E19", "C20", "E20", "E21", "E22", "C23", "E23", "C24",
    "E24", "C25", "E25", "C26", "E26", "C27", "E27", "C28", "E28", "C29",
    "E29", "E30", "E31", "C32", "E32", "C33", "E33", "C34", "E34", "E35",
    "E36", "D37", "E37", "E38", "C39", "E39", "E40", "C41", "E41", "C42",
    "E43", "C44", "C45", "E45", "C46", "E46", "C47", "E47", "C48", "D48",
    "E48", "E49", "C50", "E50"
)

foreach ($addr in $blankCells) {
    $ws.Range($addr).Value = "-"
}

# B50 did not exist as a cell at all before the edit (row 50 only had
# A/C/D/E). Give it the same number/alignment format (style index "2") as
# its row-mates before writing the value, by copying formats across instead
# of poking individual alignment properties (which would otherwise leave a
# stray, unused cellXf behind in styles.xml).
$ws.Range("A50").Copy() | Out-Null
$ws.Range("B50").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("B50").Value = "-"

# Restore the selection/active-cell state recorded in the saved workbook:
# the whole data range A2:E50 selected with A2 as the active cell.
$ws.Range("A2:E50").Select() | Out-Null
